# The commit simply restores/updates the value of cell C10 on the "Rules"
# worksheet (row 10, the R30 rule's "Integer min" value) from 18 to 1.
# All of the other differences visible in the raw OOXML diff (numFmts
# count="0", xfId="0" on cellStyleXfs, attribute re-ordering/boolean
# spelling on <col> elements, t="n" marker, "1.0" vs "1") are just
# artifacts of the resaving tool's XML serialization and not actual
# content/formatting changes, so we only need to update the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
